$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# The "Double check" text (only used by the E8:E11 "Double check"/"Fix" list,
# tied to the rule that now needs a hyphen) becomes "Double-check".
$ws.Range("E8").Value = "Double-check"
$ws.Range("E9").Value = "Double-check"
$ws.Range("E10").Value = "Double-check"
$ws.Range("E11").Value = "Double-check"

# Split the E2:E21 data validation list so that only E8:E11 offers the new
# hyphenated wording, while the rest of the column keeps the original list.
$rng = $ws.Range("E8:E11")
$rng.Validation.Delete()
$rng.Validation.Add(3, 1, 1, """Fix, Double-check""")

# Restore the default (top-left) selection/scroll position for the sheet.
$ws.Range("A1").Select()
